# New crime data collected
# Updates the NYPD 32nd Precinct weekly CompStat report:
#  - Volume/Number and "Report Covering the Week ... Through ..." header text
#  - The Crime Complaints table (rows 15-30) with refreshed weekly/28-day/
#    year-to-date/2-year figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -------------------------------------------------

# A8: "Volume 31   Number  51" -> "...Number  52"
$a8 = $ws.Range("A8")
$num = $a8.Characters(21, 2)
$num.Text = "52"
$num.Font.Name = "Andale WT"
$num.Font.Size = 10

# C9: "Report Covering the Week  12/16/2024  Through  12/22/2024"
#  -> "...12/23/2024  Through  12/29/2024"
$c9 = $ws.Range("C9")
$weekStart = $c9.Characters(27, 10)
$weekStart.Text = "12/23/2024"
$weekStart.Font.Name = "Andale WT"
$weekStart.Font.Size = 10
$weekEnd = $c9.Characters(48, 10)
$weekEnd.Text = "12/29/2024"
$weekEnd.Font.Name = "Andale WT"
$weekEnd.Font.Size = 10

# --- Crime Complaints table (rows 15-30) ----------------------------------

$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("N15").Value = -74.71264367816
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 25
$ws.Range("F16").Value = 16
$ws.Range("G16").Value = 30
$ws.Range("H16").Value = -46.666666666666
$ws.Range("I16").Value = 215
$ws.Range("J16").Value = 223
$ws.Range("K16").Value = -3.587443946188
$ws.Range("L16").Value = 0.93896713615
$ws.Range("M16").Value = -26.621160409556
$ws.Range("N16").Value = -77.368421052631
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 120
$ws.Range("G17").Value = 28
$ws.Range("H17").Value = 60.714285714285
$ws.Range("I17").Value = 492
$ws.Range("J17").Value = 449
$ws.Range("K17").Value = 9.576837416481
$ws.Range("L17").Value = -4.651162790697
$ws.Range("M17").Value = 59.223300970873
$ws.Range("N17").Value = -49.01554404145
$ws.Range("C18").Value = 2
$ws.Range("C18").NumberFormat = "#,##0"
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 8
$ws.Range("H18").Value = -11.111111111111
$ws.Range("I18").Value = 128
$ws.Range("J18").Value = 171
$ws.Range("K18").Value = -25.146198830409
$ws.Range("L18").Value = -37.560975609756
$ws.Range("M18").Value = 4.918032786885
$ws.Range("N18").Value = -86.831275720164
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 33.333333333333
$ws.Range("F19").Value = 22
$ws.Range("G19").Value = 31
$ws.Range("H19").Value = -29.032258064516
$ws.Range("I19").Value = 355
$ws.Range("J19").Value = 449
$ws.Range("K19").Value = -20.935412026726
$ws.Range("L19").Value = 4.411764705882
$ws.Range("M19").Value = 19.932432432432
$ws.Range("N19").Value = -17.249417249417
$ws.Range("D20").Value = 5
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -77.777777777777
$ws.Range("J20").Value = 139
$ws.Range("K20").Value = -41.726618705036
$ws.Range("L20").Value = 3.846153846153
$ws.Range("M20").Value = 52.830188679245
$ws.Range("N20").Value = -68.359375
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = 4.761904761904
$ws.Range("G21").Value = 110
$ws.Range("H21").Value = -10.90909090909
$ws.Range("I21").Value = 1302
$ws.Range("J21").Value = 1472
$ws.Range("K21").Value = -11.548913043478
$ws.Range("L21").Value = -5.992779783393
$ws.Range("M21").Value = 16.981132075471
$ws.Range("N21").Value = -64.952893674293
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 12
$ws.Range("K22").Value = -25
$ws.Range("L22").Value = -45.454545454545
$ws.Range("M22").Value = 20
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = -60
$ws.Range("F23").Value = 18
$ws.Range("H23").Value = -28
$ws.Range("I23").Value = 262
$ws.Range("J23").Value = 249
$ws.Range("K23").Value = 5.220883534136
$ws.Range("L23").Value = 10.548523206751
$ws.Range("M23").Value = 51.445086705202
$ws.Range("C24").Value = 15
$ws.Range("D24").Value = 6
$ws.Range("E24").Value = 150
$ws.Range("F24").Value = 76
$ws.Range("G24").Value = 64
$ws.Range("H24").Value = 18.75
$ws.Range("I24").Value = 867
$ws.Range("J24").Value = 874
$ws.Range("K24").Value = -0.800915331807
$ws.Range("L24").Value = 4.20673076923
$ws.Range("M24").Value = 36.106750392464
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 0
$ws.Range("G25").Value = 12
$ws.Range("H25").Value = 41.666666666666
$ws.Range("I25").Value = 158
$ws.Range("J25").Value = 183
$ws.Range("K25").Value = -13.661202185792
$ws.Range("L25").Value = -22.926829268292
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 44
$ws.Range("G26").Value = 47
$ws.Range("H26").Value = -6.382978723404
$ws.Range("I26").Value = 708
$ws.Range("J26").Value = 622
$ws.Range("K26").Value = 13.826366559485
$ws.Range("L26").Value = 20.203735144312
$ws.Range("M26").Value = -16.213017751479
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -100
$ws.Range("J27").Value = 49
$ws.Range("K27").Value = -36.734693877551
$ws.Range("C28").Value = 2
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 62
$ws.Range("K28").Value = 19.230769230769
$ws.Range("L28").Value = -12.676056338028
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 100
$ws.Range("L29").Value = -39.473684210526
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 100
$ws.Range("L30").Value = -37.5
